# Mark Song of Songs (row 23), Jeremiah (row 25) and Lamentations (row 26)
# as finished books ("Book done" = 1) on the "all" worksheet. The
# "Verses done" column (E) is a shared formula (IF(F=1,C,0)) and the
# totals/percentages in rows 41-42 recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")
$ws.Activate()

$ws.Range("F23").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 1

$excel.Calculate()

# Move the view back to the top of the sheet and select the whole table,
# matching the state Excel leaves the workbook in after editing near the
# bottom of the list and then pressing Ctrl+Home / selecting everything.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1:F42").Select() | Out-Null
